$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would be misread as a number by Excel's automatic
# type inference (plain decimals like "245.07") are forced to Text format
# first so the literal string from the source data is preserved verbatim.

$ws.Range("D2").Value = "30.465.51"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.912.24"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.07"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4821"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2892"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06717"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.74"
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.05"
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("D12").Value = "1.916.74"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07551"
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.260"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6709"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "287.63"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "30.465.29"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9988"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "2.163.40"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.470"
$ws.Range("E22").Value = "  +4.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9984"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.399"
$ws.Range("E24").Value = "  +3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.462"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.89"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.30"
$ws.Range("E27").Value = "  -5.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.112"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.167"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.041"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04974"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7276"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9990"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02035"
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.666"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.55"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.014"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4437"
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8660"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.800"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9986"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.13"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.335"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.06"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.259"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.77"
$ws.Range("E51").Value = "  -0.43%  "
